# "Fruta / hortaliza, semanal" — weekly refresh: a new weekly price
# observation is inserted as row 200 (pushing the former rows 200:208
# down to 201:209, and extending the used range from R208 to R209).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 200; everything at/after row 200
# (data rows + the sheet dimension) shifts down by one.
$ws.Rows.Item(200).Insert()

# Populate the newly inserted row 200 with this week's observation.
$ws.Range("A200").Value2 = 9
$ws.Range("B200").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C200").Value2 = "Metropolitana"
$ws.Range("D200").Value2 = 44568
$ws.Range("E200").Value2 = 13
$ws.Range("F200").Value2 = 300000001
$ws.Range("G200").Value2 = "Rabanito"
$ws.Range("H200").Value2 = "Sin especificar"
$ws.Range("I200").Value2 = "Primera"
$ws.Range("J200").Value2 = 6100
$ws.Range("K200").Value2 = 2500
$ws.Range("L200").Value2 = 3000
$ws.Range("M200").Value2 = 2750
$ws.Range("N200").Value2 = "`$/cien unidades (volumen en unidades)"
$ws.Range("O200").Value2 = "Provincia de Chacabuco"
$ws.Range("P200").Value2 = 28
$ws.Range("Q200").Value2 = 100
$ws.Range("R200").Value2 = "Hortaliza"
